$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row
$ws.Range("A1").Value = "ENTITY ID"
$ws.Range("B1").Value = "MIGRATION DATE"
$ws.Range("C1").Value = "ADDRESS"
$ws.Range("D1").Value = "FINANCIAL INSTITUTION NAME"

# Apply the existing bold/border/center style (already on A1:B1) to C1:D1 too
$ws.Range("A1:D1").Style = "Normal"
$ws.Range("C1").Font.Bold = $true
$ws.Range("D1").Font.Bold = $true
$ws.Range("C1:D1").Borders.LineStyle = 1
$ws.Range("C1:D1").HorizontalAlignment = -4108
$ws.Range("C1:D1").VerticalAlignment = -4160

# Data row
$ws.Range("A2").Value = "111AAA111"
$ws.Range("B2").Value = "2025-10-20"
$ws.Range("C2").Value = "Karapakkam"
$ws.Range("D2").Value = "AAA"
